# Update the "countries" (Pais) worksheet with the refreshed COVID-19 snapshot.
# - Updates the "last refreshed" timestamp in A1.
# - Updates case/recovered/critical/death counters for several countries whose
#   totals moved between snapshots.
# - Three countries (Guadalupe, Curazao) climbed past their neighbours in the
#   ranking and Nueva Caledonia/Santa Lucia swapped tie-break order, so the
#   country names shown in column A for the affected row ranges are rewritten
#   together with their statistics so the sheet stays sorted by total cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 21:14"
$ws.Range("B4").Value = 7431458
$ws.Range("C4").Value = 25312
$ws.Range("D4").Value = 4677570
$ws.Range("E4").Value = 2542495
$ws.Range("G4").Value = 608
$ws.Range("H4").Value = 211393
$ws.Range("D14").Value = 96797
$ws.Range("E14").Value = 434782
$ws.Range("B25").Value = 292559
$ws.Range("C25").Value = 2093
$ws.Range("E25").Value = 26992
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 9567
$ws.Range("B63").Value = 51530
$ws.Range("C63").Value = 162
$ws.Range("D63").Value = 36174
$ws.Range("E63").Value = 13620
$ws.Range("G63").Value = 10
$ws.Range("H63").Value = 1736
$ws.Range("B103").Value = 10291
$ws.Range("C103").Value = 97
$ws.Range("D103").Value = 9108
$ws.Range("E103").Value = 1149
$ws.Range("B117").Value = 6024
$ws.Range("C117").Value = 124
$ws.Range("D117").Value = 5277
$ws.Range("E117").Value = 687
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 60
$ws.Range("A121").Value = "Guadalupe"
$ws.Range("B121").Value = 5528
$ws.Range("C121").Value = 1041
$ws.Range("D121").Value = 2199
$ws.Range("E121").Value = 3272
$ws.Range("G121").Value = 15
$ws.Range("H121").Value = 57
$ws.Range("A122").Value = "Suazilandia"
$ws.Range("B122").Value = 5482
$ws.Range("C122").Value = 20
$ws.Range("D122").Value = 4912
$ws.Range("E122").Value = 461
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 109
$ws.Range("A123").Value = "Republica de Yibuti"
$ws.Range("B123").Value = 5416
$ws.Range("D123").Value = 5344
$ws.Range("E123").Value = 11
$ws.Range("H123").Value = 61
$ws.Range("A124").Value = "Nicaragua"
$ws.Range("B124").Value = 5170
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 2913
$ws.Range("E124").Value = 2106
$ws.Range("H124").Value = 151
$ws.Range("A125").Value = "Hong Kong"
$ws.Range("B125").Value = 5088
$ws.Range("C125").Value = 8
$ws.Range("D125").Value = 4827
$ws.Range("E125").Value = 156
$ws.Range("H125").Value = 105
$ws.Range("A126").Value = "Guinea Ecuatorial"
$ws.Range("B126").Value = 5030
$ws.Range("D126").Value = 4769
$ws.Range("E126").Value = 178
$ws.Range("H126").Value = 83
$ws.Range("A127").Value = "Congo"
$ws.Range("B127").Value = 5008
$ws.Range("D127").Value = 3887
$ws.Range("E127").Value = 1032
$ws.Range("H127").Value = 89
$ws.Range("A128").Value = "Angola"
$ws.Range("B128").Value = 4972
$ws.Range("C128").Value = 67
$ws.Range("D128").Value = 1941
$ws.Range("E128").Value = 2848
$ws.Range("G128").Value = 4
$ws.Range("H128").Value = 183
$ws.Range("A129").Value = "Surinam"
$ws.Range("B129").Value = 4863
$ws.Range("D129").Value = 4676
$ws.Range("E129").Value = 83
$ws.Range("H129").Value = 104
$ws.Range("A130").Value = "Ruanda"
$ws.Range("B130").Value = 4836
$ws.Range("D130").Value = 3125
$ws.Range("E130").Value = 1682
$ws.Range("H130").Value = 29
$ws.Range("A131").Value = "Republica de Africa Central"
$ws.Range("B131").Value = 4806
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 1840
$ws.Range("E131").Value = 2904
$ws.Range("H131").Value = 62
$ws.Range("A132").Value = "Lituania"
$ws.Range("B132").Value = 4693
$ws.Range("C132").Value = 115
$ws.Range("D132").Value = 2365
$ws.Range("E132").Value = 2236
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 92
$ws.Range("A133").Value = "Trinidad yTobago"
$ws.Range("B133").Value = 4517
$ws.Range("C133").Value = 54
$ws.Range("D133").Value = 2560
$ws.Range("E133").Value = 1882
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 75
$ws.Range("B152").Value = 2231
$ws.Range("C152").Value = 9
$ws.Range("D152").Value = 1687
$ws.Range("E152").Value = 472
$ws.Range("B177").Value = 508
$ws.Range("C177").Value = 2
$ws.Range("E177").Value = 35
$ws.Range("A181").Value = "Curazao"
$ws.Range("B181").Value = 392
$ws.Range("C181").Value = 22
$ws.Range("D181").Value = 171
$ws.Range("E181").Value = 220
$ws.Range("H181").Value = 1
$ws.Range("A182").Value = "San Martin (Parte Francesa)"
$ws.Range("B182").Value = 383
$ws.Range("D182").Value = 273
$ws.Range("E182").Value = 102
$ws.Range("H182").Value = 8
$ws.Range("A183").Value = "Eritrea"
$ws.Range("B183").Value = 375
$ws.Range("D183").Value = 341
$ws.Range("E183").Value = 34
$ws.Range("H183").Value = 0
$ws.Range("B189").Value = 218
$ws.Range("C189").Value = 4
$ws.Range("D189").Value = 181
$ws.Range("E189").Value = 35
$ws.Range("D191").Value = 179
$ws.Range("E191").Value = 4
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"
